$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.884.85"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "2.236.95"
$ws.Range("E3").Value = "  -1.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "114.80"
$ws.Range("E5").Value = "  +3.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "275.77"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  -2.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.49"
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0926"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.07"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("E13").Value = "  -3.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.27"
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.875"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("D16").Value = "2.574.70"
$ws.Range("E16").Value = "  -1.50%  "
$ws.Range("D17").Value = "2.244.16"
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("D18").Value = "42.829.23"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.12"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("E22").Value = "  -3.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.97"
$ws.Range("E23").Value = "  +4.77%  "
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.26"
$ws.Range("E25").Value = "  -1.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.92"
$ws.Range("E26").Value = "  +5.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.23"
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.27"
$ws.Range("E31").Value = "  -2.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "173.11"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  -1.45%  "
$ws.Range("E34").Value = "  -0.77%  "
$ws.Range("E35").Value = "  -1.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.41"
$ws.Range("E36").Value = "  +12.93%  "
$ws.Range("E37").Value = "  -1.75%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.64"
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0371"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.99"
$ws.Range("E42").Value = "  -6.08%  "
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.20"
$ws.Range("E44").Value = "  -7.31%  "
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.34"
$ws.Range("E46").Value = "  -2.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.64"
$ws.Range("E47").Value = "  -6.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.45"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0991"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "100.31"
$ws.Range("E51").Value = "  -0.48%  "
